$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update name and code for row 11 (previously "HAREF MACIEL" / hash)
$ws.Range("B11").Value = "MIRIAN GONÇALVES"
$ws.Range("C11").Value = "dc8f532d890d2fc187f8fdc7bf906cbf"

# Update start date and remaining days for row 11
$ws.Range("D11").Value = (Get-Date -Year 2022 -Month 9 -Day 26 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("E11").Value = 9
